$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (día, mes, año, Pruebas Realizadas, Pruebas Positivas)
$data = @(
    @(30, 11, 2020, 4926, 712),
    @(1, 12, 2020, 4529, 686),
    @(2, 12, 2020, 4074, 593),
    @(3, 12, 2020, 3776, 752),
    @(4, 12, 2020, 4039, 547),
    @(5, 12, 2020, 1307, 198),
    @(6, 12, 2020, 714, 124),
    @(7, 12, 2020, 4767, 799)
)

$startRow = 248
$endRow = $startRow + $data.Count - 1
$formula = '=+Condicion_Pacientes[[#This Row],[día]]&"/"&Condicion_Pacientes[[#This Row],[mes]]&"/"&Condicion_Pacientes[[#This Row],[año]]'

# Copy the formatting from the last existing table row down onto the new rows
# (mirrors Excel's own table auto-expand behaviour: centered day/month/year/date
# style carried forward onto freshly appended rows).
$ws.Range("A247:F247").Copy()
$ws.Range("A${startRow}:F${endRow}").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $ws.Range("B$row").Value = $rowData[0]
    $ws.Range("C$row").Value = $rowData[1]
    $ws.Range("D$row").Value = $rowData[2]
    $ws.Range("E$row").Value = $rowData[3]
    $ws.Range("F$row").Value = $rowData[4]

    $ws.Range("A$row").Formula = $formula
}

# Grow the Excel Table ("Condicion_Pacientes") to include the new rows
$tbl = $ws.ListObjects.Item("Condicion_Pacientes")
$tbl.Resize($ws.Range("A1:I$endRow"))

# Update the sheet view to match new selection/scroll position
$win = $excel.ActiveWindow
$win.ScrollRow = 239
$win.ScrollColumn = 1
$ws.Range("E256").Select() | Out-Null
